$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$helper = $ws.Range("Z1")

$ws.Range('D2').Value = '27.969.31'
$ws.Range('E2').Value = '  -0.36%  '

$ws.Range('D3').Value = '1.866.36'
$ws.Range('E3').Value = '  -1.30%  '

$helper.NumberFormat = '@'
$helper.Value = '1.006'
$helper.Copy()
$ws.Range('D4').PasteSpecial(-4163)
$ws.Range('E4').Value = '  +0.46%  '

$helper.NumberFormat = '@'
$helper.Value = '311.96'
$helper.Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  -0.64%  '

$ws.Range('E6').Value = '  +0.22%  '

$helper.NumberFormat = '@'
$helper.Value = '0.5067'
$helper.Copy()
$ws.Range('D7').PasteSpecial(-4163)
$ws.Range('E7').Value = '  +1.12%  '

$helper.NumberFormat = '@'
$helper.Value = '0.3800'
$helper.Copy()
$ws.Range('D8').PasteSpecial(-4163)
$ws.Range('E8').Value = '  -2.25%  '

$helper.NumberFormat = '@'
$helper.Value = '0.08277'
$helper.Copy()
$ws.Range('D9').PasteSpecial(-4163)
$ws.Range('E9').Value = '  -9.81%  '

$helper.NumberFormat = '@'
$helper.Value = '1.106'
$helper.Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('E10').Value = '  -1.98%  '

$helper.NumberFormat = '@'
$helper.Value = '41.40'
$helper.Copy()
$ws.Range('D11').PasteSpecial(-4163)

$helper.NumberFormat = '@'
$helper.Value = '6.197'
$helper.Copy()
$ws.Range('D12').PasteSpecial(-4163)
$ws.Range('E12').Value = '  -2.74%  '

$ws.Range('D13').Value = '1.864.04'
$ws.Range('E13').Value = '  -1.26%  '

$helper.NumberFormat = '@'
$helper.Value = '20.39'
$helper.Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('E14').Value = '  -2.07%  '

$helper.NumberFormat = '@'
$helper.Value = '7.168'
$helper.Copy()
$ws.Range('D15').PasteSpecial(-4163)
$ws.Range('E15').Value = '  -1.91%  '

$helper.NumberFormat = '@'
$helper.Value = '1.004'
$helper.Copy()
$ws.Range('D16').PasteSpecial(-4163)
$ws.Range('E16').Value = '  +0.20%  '

$helper.NumberFormat = '@'
$helper.Value = '0.00001091'
$helper.Copy()
$ws.Range('D17').PasteSpecial(-4163)
$ws.Range('E17').Value = '  -1.36%  '

$helper.NumberFormat = '@'
$helper.Value = '90.34'
$helper.Copy()
$ws.Range('D18').PasteSpecial(-4163)
$ws.Range('E18').Value = '  -1.61%  '

$helper.NumberFormat = '@'
$helper.Value = '0.06617'
$helper.Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Value = '  -0.23%  '

$helper.NumberFormat = '@'
$helper.Value = '17.80'
$helper.Copy()
$ws.Range('D20').PasteSpecial(-4163)
$ws.Range('E20').Value = '  -0.48%  '

$helper.NumberFormat = '@'
$helper.Value = '5.988'
$helper.Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('E22').Value = '  -3.94%  '

$ws.Range('D23').Value = '28.041.19'
$ws.Range('E23').Value = '  -0.34%  '

$helper.NumberFormat = '@'
$helper.Value = '11.08'
$helper.Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Value = '  -2.62%  '

$helper.NumberFormat = '@'
$helper.Value = '2.256'
$helper.Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  -2.66%  '

$helper.NumberFormat = '@'
$helper.Value = '2.553'
$helper.Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('E26').Value = '  -0.08%  '

$ws.Range('E27').Value = '  -1.01%  '

$helper.NumberFormat = '@'
$helper.Value = '157.41'
$helper.Copy()
$ws.Range('D28').PasteSpecial(-4163)
$ws.Range('E28').Value = '  -0.78%  '

$helper.NumberFormat = '@'
$helper.Value = '20.40'
$helper.Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Value = '  -1.74%  '

$helper.NumberFormat = '@'
$helper.Value = '125.34'
$helper.Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('E30').Value = '  -1.29%  '

$helper.NumberFormat = '@'
$helper.Value = '0.1054'
$helper.Copy()
$ws.Range('D31').PasteSpecial(-4163)
$ws.Range('E31').Value = '  +0.07%  '

$helper.NumberFormat = '@'
$helper.Value = '1.035'
$helper.Copy()
$ws.Range('D32').PasteSpecial(-4163)
$ws.Range('E32').Value = '  -3.46%  '

$helper.NumberFormat = '@'
$helper.Value = '5.577'
$helper.Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('E33').Value = '  -0.35%  '

$helper.NumberFormat = '@'
$helper.Value = '3.595'
$helper.Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('E34').Value = '  +0.08%  '

$helper.NumberFormat = '@'
$helper.Value = '9.629'
$helper.Copy()
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('E35').Value = '  +1.95%  '

$helper.NumberFormat = '@'
$helper.Value = '0.02423'
$helper.Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('E36').Value = '  +0.30%  '

$helper.NumberFormat = '@'
$helper.Value = '0.06516'
$helper.Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('E37').Value = '  -0.94%  '

$helper.NumberFormat = '@'
$helper.Value = '0.2151'
$helper.Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('E38').Value = '  -1.87%  '

$helper.NumberFormat = '@'
$helper.Value = '1.204'
$helper.Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Value = '  -0.76%  '

$helper.NumberFormat = '@'
$helper.Value = '0.6412'
$helper.Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Value = '  -0.20%  '

$helper.NumberFormat = '@'
$helper.Value = '1.236'
$helper.Copy()
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('E41').Value = '  -7.35%  '

$helper.NumberFormat = '@'
$helper.Value = '11.22'
$helper.Copy()
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('E42').Value = '  -2.95%  '

$helper.NumberFormat = '@'
$helper.Value = '4.852'
$helper.Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('E43').Value = '  -1.93%  '

$helper.NumberFormat = '@'
$helper.Value = '0.6072'
$helper.Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  +0.45%  '

$helper.NumberFormat = '@'
$helper.Value = '13.05'
$helper.Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Value = '  -2.26%  '

$helper.NumberFormat = '@'
$helper.Value = '1.286'
$helper.Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  -1.16%  '

$helper.NumberFormat = '@'
$helper.Value = '3.660'
$helper.Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('E47').Value = '  -0.68%  '

$helper.NumberFormat = '@'
$helper.Value = '1.992'
$helper.Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Value = '  -0.40%  '

$helper.NumberFormat = '@'
$helper.Value = '1.208'
$helper.Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('E49').Value = '  +0.52%  '

$helper.NumberFormat = '@'
$helper.Value = '121.09'
$helper.Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('E50').Value = '  -0.16%  '

$helper.NumberFormat = '@'
$helper.Value = '79.58'
$helper.Copy()
$ws.Range('D51').PasteSpecial(-4163)
$ws.Range('E51').Value = '  +0.58%  '

$excel.CutCopyMode = $false
$helper.Clear()